$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The grader rubric's milestone markers that read "III" were re-marked to "I"
# (all occurrences, so the shared string itself is retargeted).
$ws.Range("E4").Value = "I"
$ws.Range("E17").Value = "I"
$ws.Range("E60").Value = "I"
$ws.Range("E61").Value = "I"

# Clear the stray "X" marks that had been copy/pasted into D83:E84 next to C83:C84.
$ws.Range("D83").ClearContents()
$ws.Range("E83").ClearContents()
$ws.Range("D84").ClearContents()
$ws.Range("E84").ClearContents()

# Update the saved selection state to match what was active when the file was saved.
$ws.Range("F5").Select()
